$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (timestamp) column for rows 2-10 to reflect the new run time
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-12 18:28:15"
}

# Update the title text for row 10
$ws.Cells.Item(10, 2).Value = "微生物の特定と分類を行いたく、画像解析の専門家を探しています!(急いでません!)"
